$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Positions")

# Row 2 (AAPL) - updated market price / value snapshot
$ws.Range("G2").Value = 195.6349945068359
$ws.Range("H2").Value = 3912.699890136719
$ws.Range("I2").Value = 1301.699890136719

# Row 3 (GME) - updated market price / value snapshot
$ws.Range("G3").Value = 23.28000068664551
$ws.Range("H3").Value = 232.8000068664551
$ws.Range("I3").Value = -769.1999931335449
